$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (changed) date column C for all data rows (2-28) to new value
$ws.Range("C2:C28").Value = 45488

# Remove the last data row (row 29) entirely - shifts dimension/rows up
$ws.Rows.Item(29).Delete()

# Row 28 loses its explicit custom row height after the delete; autofit clears
# the leftover ht/customHeight attributes so it matches the default row.
$ws.Rows.Item(28).AutoFit()
